# Refresh market-price-derived columns (H-N) across the crafting-class sheets.
# Values correspond to a scheduled market data refresh; row/column layout is unchanged
# except for a few rows where a LeveProfitNQ (M) or LeveProfitHQ (N) cell is
# added/removed because the underlying NQ/HQ price became zero/non-zero.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 82
$ws.Cells.Item(82, 8).Value = 6055.8667
$ws.Cells.Item(82, 9).Value = 1742.625
$ws.Cells.Item(82, 11).Value = 5227.875
$ws.Cells.Item(82, 13).Value = -4821.875

# Row 85
$ws.Cells.Item(85, 8).Value = 6055.8667
$ws.Cells.Item(85, 9).Value = 1742.625
$ws.Cells.Item(85, 11).Value = 5227.875
$ws.Cells.Item(85, 13).Value = -3823.875

# Row 113
$ws.Cells.Item(113, 8).Value = 8936.429
$ws.Cells.Item(113, 9).Value = 4029.1428
$ws.Cells.Item(113, 10).Value = 13843.714
$ws.Cells.Item(113, 11).Value = 4029.1428
$ws.Cells.Item(113, 12).Value = 13843.714
$ws.Cells.Item(113, 13).Value = -775.1428000000001
$ws.Cells.Item(113, 14).Value = -20351.714

# Row 125
$ws.Cells.Item(125, 8).Value = 4059.2307
$ws.Cells.Item(125, 9).Value = 3238
$ws.Cells.Item(125, 10).Value = 4572.5
$ws.Cells.Item(125, 11).Value = 29142
$ws.Cells.Item(125, 12).Value = 41152.5
$ws.Cells.Item(125, 13).Value = -26682
$ws.Cells.Item(125, 14).Value = -46072.5

# Row 137
$ws.Cells.Item(137, 8).Value = 3158.3076
$ws.Cells.Item(137, 9).Value = 2096.182
$ws.Cells.Item(137, 11).Value = 6288.545999999999
$ws.Cells.Item(137, 13).Value = -3738.545999999999

# Row 141
$ws.Cells.Item(141, 8).Value = 68919.39999999999
$ws.Cells.Item(141, 9).Value = 78914.69500000001
$ws.Cells.Item(141, 10).Value = 3950
$ws.Cells.Item(141, 11).Value = 236744.085
$ws.Cells.Item(141, 12).Value = 11850
$ws.Cells.Item(141, 13).Value = -231564.085
$ws.Cells.Item(141, 14).Value = -22210

$ws = $wb.Worksheets.Item("ARM")
# Row 27
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 14).ClearContents()  # N27 removed

# Row 32
$ws.Cells.Item(32, 8).Value = 6409.6416
$ws.Cells.Item(32, 9).Value = 4800.196
$ws.Cells.Item(32, 10).Value = 9935.096
$ws.Cells.Item(32, 11).Value = 4800.196
$ws.Cells.Item(32, 12).Value = 9935.096
$ws.Cells.Item(32, 13).Value = -4513.196
$ws.Cells.Item(32, 14).Value = -10509.096

# Row 48
$ws.Cells.Item(48, 8).Value = 79800
$ws.Cells.Item(48, 10).Value = 79800
$ws.Cells.Item(48, 12).Value = 79800
$ws.Cells.Item(48, 14).Value = -80568

# Row 122
$ws.Cells.Item(122, 8).Value = 1914.6364
$ws.Cells.Item(122, 9).Value = 1362.8889
$ws.Cells.Item(122, 10).Value = 4397.5
$ws.Cells.Item(122, 11).Value = 4088.6667
$ws.Cells.Item(122, 12).Value = 13192.5
$ws.Cells.Item(122, 13).Value = -1638.6667
$ws.Cells.Item(122, 14).Value = -18092.5

# Row 132
$ws.Cells.Item(132, 8).Value = 2173.6604
$ws.Cells.Item(132, 9).Value = 1676.9025
$ws.Cells.Item(132, 10).Value = 3870.9167
$ws.Cells.Item(132, 11).Value = 5030.7075
$ws.Cells.Item(132, 12).Value = 11612.7501
$ws.Cells.Item(132, 13).Value = -2500.7075
$ws.Cells.Item(132, 14).Value = -16672.7501

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Cells.Item(134, 8).Value = 3070.0688
$ws.Cells.Item(134, 9).Value = 1321.2
$ws.Cells.Item(134, 10).Value = 4943.857
$ws.Cells.Item(134, 11).Value = 3963.6
$ws.Cells.Item(134, 12).Value = 14831.571
$ws.Cells.Item(134, 13).Value = -1428.6
$ws.Cells.Item(134, 14).Value = -19901.571

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 11113920
$ws.Cells.Item(31, 9).Value = 1385.4667
$ws.Cells.Item(31, 11).Value = 1385.4667
$ws.Cells.Item(31, 13).Value = -1090.4667

# Row 34
$ws.Cells.Item(34, 8).Value = 11113920
$ws.Cells.Item(34, 9).Value = 1385.4667
$ws.Cells.Item(34, 11).Value = 1385.4667
$ws.Cells.Item(34, 13).Value = -1183.4667

# Row 134
$ws.Cells.Item(134, 8).Value = 6528.56
$ws.Cells.Item(134, 10).Value = 3455.7144
$ws.Cells.Item(134, 12).Value = 10367.1432
$ws.Cells.Item(134, 14).Value = -15437.1432

# Row 135
$ws.Cells.Item(135, 8).Value = 37281.668
$ws.Cells.Item(135, 10).Value = 37281.668
$ws.Cells.Item(135, 12).Value = 37281.668
$ws.Cells.Item(135, 14).Value = -47421.668

$ws = $wb.Worksheets.Item("CUL")
# Row 50
$ws.Cells.Item(50, 8).Value = 520.6667
$ws.Cells.Item(50, 9).Value = 374.5
$ws.Cells.Item(50, 10).Value = 573.8182
$ws.Cells.Item(50, 11).Value = 1123.5
$ws.Cells.Item(50, 12).Value = 1721.4546
$ws.Cells.Item(50, 13).Value = -642.5
$ws.Cells.Item(50, 14).Value = -2683.4546

# Row 53
$ws.Cells.Item(53, 8).Value = 520.6667
$ws.Cells.Item(53, 9).Value = 374.5
$ws.Cells.Item(53, 10).Value = 573.8182
$ws.Cells.Item(53, 11).Value = 1123.5
$ws.Cells.Item(53, 12).Value = 1721.4546
$ws.Cells.Item(53, 13).Value = -642.5
$ws.Cells.Item(53, 14).Value = -2683.4546

# Row 55
$ws.Cells.Item(55, 8).Value = 1250
$ws.Cells.Item(55, 9).Value = 500
$ws.Cells.Item(55, 10).Value = 2000
$ws.Cells.Item(55, 11).Value = 1500
$ws.Cells.Item(55, 12).Value = 6000
$ws.Cells.Item(55, 13).Value = -1323
$ws.Cells.Item(55, 14).Value = -6354

# Row 86
$ws.Cells.Item(86, 8).Value = 951
$ws.Cells.Item(86, 9).Value = 951
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 2853
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = -1667
$ws.Cells.Item(86, 14).ClearContents()  # N86 removed

# Row 89
$ws.Cells.Item(89, 8).Value = 951
$ws.Cells.Item(89, 9).Value = 951
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 8559
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).Value = -2631
$ws.Cells.Item(89, 14).ClearContents()  # N89 removed

# Row 92
$ws.Cells.Item(92, 8).Value = 3703.8333
$ws.Cells.Item(92, 10).Value = 5313.25
$ws.Cells.Item(92, 12).Value = 15939.75
$ws.Cells.Item(92, 14).Value = -18435.75

# Row 131
$ws.Cells.Item(131, 8).Value = 8624951
$ws.Cells.Item(131, 9).Value = 125050200
$ws.Cells.Item(131, 10).Value = 858.0741
$ws.Cells.Item(131, 11).Value = 375150600
$ws.Cells.Item(131, 12).Value = 2574.2223
$ws.Cells.Item(131, 13).Value = -375145560
$ws.Cells.Item(131, 14).Value = -12654.2223

# Row 137
$ws.Cells.Item(137, 8).Value = 4151.778
$ws.Cells.Item(137, 9).Value = 666.6667
$ws.Cells.Item(137, 10).Value = 5894.3335
$ws.Cells.Item(137, 11).Value = 2000.0001
$ws.Cells.Item(137, 12).Value = 17683.0005
$ws.Cells.Item(137, 13).Value = 3099.9999
$ws.Cells.Item(137, 14).Value = -27883.0005

$ws = $wb.Worksheets.Item("GSM")
# Row 124
$ws.Cells.Item(124, 8).Value = 41722.223
$ws.Cells.Item(124, 10).Value = 41722.223
$ws.Cells.Item(124, 12).Value = 41722.223
$ws.Cells.Item(124, 14).Value = -51542.223

# Row 132
$ws.Cells.Item(132, 8).Value = 2950.027
$ws.Cells.Item(132, 9).Value = 2028.3
$ws.Cells.Item(132, 10).Value = 4034.4119
$ws.Cells.Item(132, 11).Value = 6084.9
$ws.Cells.Item(132, 12).Value = 12103.2357
$ws.Cells.Item(132, 13).Value = -3554.9
$ws.Cells.Item(132, 14).Value = -17163.2357

# Row 136
$ws.Cells.Item(136, 8).Value = 23871.875
$ws.Cells.Item(136, 10).Value = 23871.875
$ws.Cells.Item(136, 12).Value = 71615.625
$ws.Cells.Item(136, 14).Value = -76715.625

# Row 137
$ws.Cells.Item(137, 8).Value = 42661.668
$ws.Cells.Item(137, 10).Value = 42661.668
$ws.Cells.Item(137, 12).Value = 42661.668
$ws.Cells.Item(137, 14).Value = -52861.668

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Cells.Item(122, 8).Value = 4662.381
$ws.Cells.Item(122, 10).Value = 8016.2856
$ws.Cells.Item(122, 12).Value = 24048.8568
$ws.Cells.Item(122, 14).Value = -28948.8568

# Row 136
$ws.Cells.Item(136, 8).Value = 3744.25
$ws.Cells.Item(136, 9).Value = 1756.5883
$ws.Cells.Item(136, 11).Value = 5269.7649
$ws.Cells.Item(136, 13).Value = -2719.7649

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Cells.Item(136, 8).Value = 1975
$ws.Cells.Item(136, 9).Value = 864.1667
$ws.Cells.Item(136, 10).Value = 4641
$ws.Cells.Item(136, 11).Value = 2592.5001
$ws.Cells.Item(136, 12).Value = 13923
$ws.Cells.Item(136, 13).Value = -42.5001000000002
$ws.Cells.Item(136, 14).Value = -19023
